# SyntheticDataPipeline: don't generate Work dcterms:created, dcterms:date, and
# other date predicates; generate WorkClosing and WorkOpening in addition to
# WorkCreation events.

$wb = $excel.ActiveWorkbook

# --- 1. Sheets: split the single "WorkCreation" event sheet into three event
#        sheets - WorkClosing, WorkCreation, WorkOpening - right after "Work".
#        Renaming the existing sheet keeps its sheetId/content ("@graph" in
#        A1); the two new sheets are made via Copy() (rather than Add()) so
#        they inherit the same sheet formatting (outline/page setup props)
#        as the original instead of engine defaults, and their "@graph" A1
#        cell comes along for free. They pick up the next sheetIds in
#        sequence, matching the target workbook layout exactly.
$workClosing = $wb.Worksheets.Item("WorkCreation")
$workClosing.Name = "WorkClosing"

$workClosing.Copy([System.Type]::Missing, $workClosing)
$newWorkCreation = $wb.Worksheets.Item("WorkClosing (2)")
$newWorkCreation.Name = "WorkCreation"

$newWorkCreation.Copy([System.Type]::Missing, $newWorkCreation)
$newWorkOpening = $wb.Worksheets.Item("WorkCreation (2)")
$newWorkOpening.Name = "WorkOpening"

# --- 2. License sheet: add four new leading columns (cc:legalcode,
#        cc:licenseClass, cc:permits, cc:requires) between @id and the
#        existing identifier/title columns.
$license = $wb.Worksheets.Item("License")
$license.Range("B1:E1").EntireColumn.Insert()
$license.Range("B1").Value = "cc:legalcode"
$license.Range("C1").Value = "cc:licenseClass"
$license.Range("D1").Value = "cc:permits"
$license.Range("E1").Value = "cc:requires"

# title corrections that came along with the license table rebuild
$license.Range("G2").Value = "BSD License"
$license.Range("G38").Value = "MIT"

# --- 3. RightsStatement sheet: note (column E) values re-shuffled across a
#        handful of rows.
$rightsStatement = $wb.Worksheets.Item("RightsStatement")
$rightsStatement.Range("E3").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."
$rightsStatement.Range("E5").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."
$rightsStatement.Range("E6").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$rightsStatement.Range("E8").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."
$rightsStatement.Range("E9").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."
$rightsStatement.Range("E11").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$rightsStatement.Range("E13").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."

# --- 4. Person sheet: relation (column F) now points at Wikidata instead of
#        Wikipedia for rows 2, 4 and 6.
$person = $wb.Worksheets.Item("Person")
$person.Range("F2").Value = "http://www.wikidata.org/entity/Q7251"
$person.Range("F4").Value = "http://www.wikidata.org/entity/Q7251"
$person.Range("F6").Value = "http://www.wikidata.org/entity/Q7251"
